# "aggiornamento 15, 16, 17 marzo" - append three new daily rows
# (227-229) to the data table, continuing the existing series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (date number format, border, alignment = style "2")
# from the last existing data row (226) onto the three new rows before
# writing values, so the new date cells pick up the same look as the rest
# of column A.
$ws.Range("A226").Copy()
$ws.Range("A227:A229").PasteSpecial(-4122)

$ws.Range("A227").Value = 44301
$ws.Range("B227").Value = 1
$ws.Range("C227").Value = 8
$ws.Range("D227").Value = 133.0893362169356

$ws.Range("A228").Value = 44302
$ws.Range("B228").Value = 0
$ws.Range("C228").Value = 7
$ws.Range("D228").Value = 116.4531691898187

$ws.Range("A229").Value = 44303
$ws.Range("B229").Value = 1
$ws.Range("C229").Value = 8
$ws.Range("D229").Value = 133.0893362169356
